# Daily attendance processing - 2025-12-25 23:31:52
# Normalizes the "Recorded By" (column G) values: when the value starts
# with "System, ", move that leading "System" token to the end of the
# comma-separated list of recorders.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value()

    if ($val -ne $null -and $val -is [string] -and $val.StartsWith("System, ")) {
        $rest = $val.Substring(8)
        $newVal = $rest + ", System"
        $cell.Value = $newVal
    }
}
